# Updates the cryptos list with freshly scraped prices / volume percentages.
# Mirrors the data refresh performed by the "Updated cryptos list" GitHub Action.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for numeric-looking price strings so that Excel
# does not silently coerce them into floating point numbers (which would
# both lose trailing zeros and introduce binary floating point noise).
$textCells = @(
    "D5", "D6", "D9", "D10", "D11", "D12", "D14", "D15",
    "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27",
    "D28", "D30", "D31", "D34", "D35", "D36", "D37", "D38",
    "D39", "D40", "D41", "D43", "D44", "D46", "D48", "D49"
)
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Coin / Link / Price / Volume(1h) values to write.
$updates = @(
    @{ Row = 2; D = '65.870.98'; E = '  +0.56%  ' },
    @{ Row = 3; D = '3.604.96'; E = '  +1.71%  ' },
    @{ Row = 4; E = '  -0.03%  ' },
    @{ Row = 5; D = '605.46'; E = '  +0.25%  ' },
    @{ Row = 6; D = '136.67'; E = '  -2.90%  ' },
    @{ Row = 7; D = '3.604.65'; E = '  +1.72%  ' },
    @{ Row = 8; E = '  +0.06%  ' },
    @{ Row = 9; D = '0.498'; E = '  +0.81%  ' },
    @{ Row = 10; D = '0.126'; E = '  -0.08%  ' },
    @{ Row = 11; D = '7.23'; E = '  +2.62%  ' },
    @{ Row = 12; D = '0.393'; E = '  -0.45%  ' },
    @{ Row = 13; D = '4.217.22'; E = '  +1.60%  ' },
    @{ Row = 14; D = '28.16'; E = '  +3.21%  ' },
    @{ Row = 15; D = '0.0000187'; E = '  -0.47%  ' },
    @{ Row = 16; D = '3.591.31'; E = '  +1.02%  ' },
    @{ Row = 18; D = '65.070.79'; E = '  -0.67%  ' },
    @{ Row = 19; D = '10.12'; E = '  -1.70%  ' },
    @{ Row = 20; D = '14.67'; E = '  +2.39%  ' },
    @{ Row = 21; D = '5.95'; E = '  +0.06%  ' },
    @{ Row = 22; D = '396.95'; E = '  +0.18%  ' },
    @{ Row = 23; D = '0.590'; E = '  +2.87%  ' },
    @{ Row = 24; D = '3.747.32'; E = '  +1.43%  ' },
    @{ Row = 25; D = '74.72'; E = '  +0.60%  ' },
    @{ Row = 26; D = '0.998'; E = '  -0.18%  ' },
    @{ Row = 27; D = '0.0000118'; E = '  +0.59%  ' },
    @{ Row = 28; D = '8.13'; E = '  +2.49%  ' },
    @{ Row = 29; E = '  +27.77%  ' },
    @{ Row = 30; D = '2.40'; E = '  +4.02%  ' },
    @{ Row = 31; D = '8.73'; E = '  +4.55%  ' },
    @{ Row = 32; E = '  -0.06%  ' },
    @{ Row = 33; D = '3.601.03'; E = '  +1.05%  ' },
    @{ Row = 34; D = '24.62'; E = '  +3.37%  ' },
    @{ Row = 35; B = 'Kaspa'; C = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D = '0.148'; E = '  +0.15%  ' },
    @{ Row = 36; B = 'USDe'; C = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'; D = '1.00'; E = '  +0.00%  ' },
    @{ Row = 37; D = '5.33'; E = '  +5.56%  ' },
    @{ Row = 38; D = '7.11'; E = '  +0.77%  ' },
    @{ Row = 39; D = '1.60'; E = '  +2.49%  ' },
    @{ Row = 40; D = '171.67'; E = '  +1.43%  ' },
    @{ Row = 41; D = '0.0834'; E = '  +2.02%  ' },
    @{ Row = 42; E = '  +1.16%  ' },
    @{ Row = 43; D = '26.38'; E = '  +0.19%  ' },
    @{ Row = 44; D = '43.39'; E = '  +0.97%  ' },
    @{ Row = 45; E = '  +2.46%  ' },
    @{ Row = 46; D = '4.54'; E = '  +2.01%  ' },
    @{ Row = 47; E = '  -0.04%  ' },
    @{ Row = 48; D = '1.71'; E = '  +0.14%  ' },
    @{ Row = 49; D = '7.09'; E = '  +3.73%  ' },
    @{ Row = 50; D = '2.467.93'; E = '  -0.17%  ' },
    @{ Row = 51; E = '  +2.94%  ' }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Range("B$($u.Row)").Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C$($u.Row)").Value = $u.C }
    if ($u.ContainsKey("D")) { $ws.Range("D$($u.Row)").Value = $u.D }
    if ($u.ContainsKey("E")) { $ws.Range("E$($u.Row)").Value = $u.E }
}
